$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 1109810.1
$ws.Range("I70").Value = 3049477.8
$ws.Range("K70").Value = 9148433.399999999
$ws.Range("M70").Value = -9148163.399999999

$ws.Range("H73").Value = 1109810.1
$ws.Range("I73").Value = 3049477.8
$ws.Range("K73").Value = 9148433.399999999
$ws.Range("M73").Value = -9147497.399999999

$ws.Range("H135").Value = 2518.8
$ws.Range("I135").Value = 1942.6666
$ws.Range("K135").Value = 17483.9994
$ws.Range("M135").Value = -14948.9994

$ws.Range("H137").Value = 1751.8235
$ws.Range("I137").Value = 1550.9546
$ws.Range("J137").Value = 2120.0833
$ws.Range("K137").Value = 4652.8638
$ws.Range("L137").Value = 6360.249899999999
$ws.Range("M137").Value = -2102.8638
$ws.Range("N137").Value = -11460.2499

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H30").Value = 40244.5
$ws.Range("I30").Value = 500
$ws.Range("K30").Value = 500
$ws.Range("M30").Value = -350

$ws.Range("H61").Value = 23207874
$ws.Range("I61").Value = 29443578
$ws.Range("J61").Value = 2006479.6
$ws.Range("K61").Value = 29443578
$ws.Range("L61").Value = 2006479.6
$ws.Range("M61").Value = -29443366
$ws.Range("N61").Value = -2006903.6

$ws.Range("H74").Value = 1745.7368
$ws.Range("I74").Value = 1534.6
$ws.Range("K74").Value = 1534.6
$ws.Range("M74").Value = -660.5999999999999

$ws.Range("H77").Value = 1745.7368
$ws.Range("I77").Value = 1534.6
$ws.Range("K77").Value = 7673
$ws.Range("M77").Value = -3305

$ws.Range("H102").Value = 22729370
$ws.Range("I102").Value = 31251708
$ws.Range("K102").Value = 31251708
$ws.Range("M102").Value = -31250086

$ws.Range("H122").Value = 6629.6875
$ws.Range("I122").Value = 5071.7334
$ws.Range("K122").Value = 15215.2002
$ws.Range("M122").Value = -12765.2002

$ws.Range("H136").Value = 23207874
$ws.Range("I136").Value = 29443578
$ws.Range("J136").Value = 2006479.6
$ws.Range("K136").Value = 88330734
$ws.Range("L136").Value = 6019438.800000001
$ws.Range("M136").Value = -88328184
$ws.Range("N136").Value = -6024538.800000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 2797.6875
$ws.Range("I80").Value = 1686.8572
$ws.Range("J80").Value = 3661.6667
$ws.Range("K80").Value = 1686.8572
$ws.Range("L80").Value = 3661.6667
$ws.Range("M80").Value = -688.8571999999999
$ws.Range("N80").Value = -5657.6667

$ws.Range("H83").Value = 2797.6875
$ws.Range("I83").Value = 1686.8572
$ws.Range("J83").Value = 3661.6667
$ws.Range("K83").Value = 8434.286
$ws.Range("L83").Value = 18308.3335
$ws.Range("M83").Value = -3442.286
$ws.Range("N83").Value = -28292.3335

$ws.Range("H94").Value = 1888.75
$ws.Range("I94").Value = 2060.0454
$ws.Range("K94").Value = 2060.0454
$ws.Range("M94").Value = -1609.0454

$ws.Range("H99").Value = 2701.4707
$ws.Range("I99").Value = 2361.6
$ws.Range("K99").Value = 2361.6
$ws.Range("M99").Value = -863.5999999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 788.0833
$ws.Range("I22").Value = 676.1429000000001
$ws.Range("J22").Value = 944.8
$ws.Range("K22").Value = 676.1429000000001
$ws.Range("L22").Value = 944.8
$ws.Range("M22").Value = -326.1429000000001
$ws.Range("N22").Value = -1644.8

$ws.Range("H28").Value = 42213.43
$ws.Range("J28").Value = 42213.43
$ws.Range("L28").Value = 42213.43
$ws.Range("N28").Value = -42703.43

$ws.Range("H70").Value = 99999
$ws.Range("J70").Value = 99999
$ws.Range("L70").Value = 99999
$ws.Range("N70").Value = -100629

$ws.Range("H73").Value = 99999
$ws.Range("J73").Value = 99999
$ws.Range("L73").Value = 99999
$ws.Range("N73").Value = -102183

$ws.Range("H74").Value = 73331.664
$ws.Range("J74").Value = 73331.664
$ws.Range("L74").Value = 73331.664
$ws.Range("N74").Value = -75079.664

$ws.Range("H77").Value = 73331.664
$ws.Range("J77").Value = 73331.664
$ws.Range("L77").Value = 219994.992
$ws.Range("N77").Value = -228730.992

$ws.Range("H81").Value = 99174.336
$ws.Range("J81").Value = 99174.336
$ws.Range("L81").Value = 99174.336
$ws.Range("N81").Value = -101170.336

$ws.Range("H82").Value = 68610.71000000001
$ws.Range("I82").Value = 60000
$ws.Range("J82").Value = 70045.836
$ws.Range("K82").Value = 60000
$ws.Range("L82").Value = 70045.836
$ws.Range("M82").Value = -59639
$ws.Range("N82").Value = -70767.836

$ws.Range("H84").Value = 99174.336
$ws.Range("J84").Value = 99174.336
$ws.Range("L84").Value = 297523.008
$ws.Range("N84").Value = -307507.008

$ws.Range("H85").Value = 68610.71000000001
$ws.Range("I85").Value = 60000
$ws.Range("J85").Value = 70045.836
$ws.Range("K85").Value = 60000
$ws.Range("L85").Value = 70045.836
$ws.Range("M85").Value = -58752
$ws.Range("N85").Value = -72541.836

$ws.Range("H88").Value = 25411.154
$ws.Range("J88").Value = 25411.154
$ws.Range("L88").Value = 25411.154
$ws.Range("N88").Value = -26223.154

$ws.Range("H91").Value = 25411.154
$ws.Range("J91").Value = 25411.154
$ws.Range("L91").Value = 25411.154
$ws.Range("N91").Value = -28219.154

$ws.Range("H122").Value = 2781.724
$ws.Range("I122").Value = 2515.4211
$ws.Range("J122").Value = 3287.7
$ws.Range("K122").Value = 7546.263300000001
$ws.Range("L122").Value = 9863.099999999999
$ws.Range("M122").Value = -5096.263300000001
$ws.Range("N122").Value = -14763.1

$ws.Range("H132").Value = 4324.75
$ws.Range("I132").Value = 4324.75
$ws.Range("K132").Value = 12974.25
$ws.Range("M132").Value = -10444.25

$ws.Range("H134").Value = 3004.0571
$ws.Range("I134").Value = 3134.2424
$ws.Range("J134").Value = 856
$ws.Range("K134").Value = 9402.727200000001
$ws.Range("L134").Value = 2568
$ws.Range("M134").Value = -6867.727200000001
$ws.Range("N134").Value = -7638

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 6405085
$ws.Range("I4").Value = 5891526
$ws.Range("K4").Value = 17674578
$ws.Range("M4").Value = -17674466

$ws.Range("H18").Value = 538.6667
$ws.Range("I18").Value = 481
$ws.Range("K18").Value = 1443
$ws.Range("M18").Value = -1274

$ws.Range("H118").Value = 7829.4546
$ws.Range("I118").Value = 5279.1
$ws.Range("K118").Value = 15837.3
$ws.Range("M118").Value = -14594.3

$ws.Range("H140").Value = 3935.2727
$ws.Range("I140").Value = 1729.125
$ws.Range("K140").Value = 5187.375
$ws.Range("M140").Value = -7.375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2074.3076
$ws.Range("I122").Value = 2246.5
$ws.Range("K122").Value = 6739.5
$ws.Range("M122").Value = -4289.5

$ws.Range("H134").Value = 82170.836
$ws.Range("J134").Value = 82170.836
$ws.Range("L134").Value = 246512.508
$ws.Range("N134").Value = -251582.508

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2061006.1
$ws.Range("I93").Value = 1628
$ws.Range("J93").Value = 4635229
$ws.Range("K93").Value = 1628
$ws.Range("L93").Value = 4635229
$ws.Range("M93").Value = -380
$ws.Range("N93").Value = -4637725

$ws.Range("H115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").ClearContents()

$ws.Range("H135").Value = 64000
$ws.Range("I135").Value = 40000
$ws.Range("K135").Value = 40000
$ws.Range("M135").Value = -34930

$ws.Range("H139").Value = 150715
$ws.Range("J139").Value = 150715
$ws.Range("L139").Value = 150715
$ws.Range("N139").Value = -160995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 7107.615
$ws.Range("J96").Value = 10113.286
$ws.Range("L96").Value = 10113.286
$ws.Range("N96").Value = -12859.286

$ws.Range("H100").Value = 657.7059
$ws.Range("I100").Value = 686.1667
$ws.Range("K100").Value = 1372.3334
$ws.Range("M100").Value = -831.3334

$ws.Range("H132").Value = 4597.5
$ws.Range("J132").Value = 3995
$ws.Range("L132").Value = 11985
$ws.Range("N132").Value = -17045

$ws.Range("H136").Value = 6095.5635
$ws.Range("I136").Value = 6466.8477
$ws.Range("J136").Value = 4197.8887
$ws.Range("K136").Value = 19400.5431
$ws.Range("L136").Value = 12593.6661
$ws.Range("M136").Value = -16850.5431
$ws.Range("N136").Value = -17693.6661

$ws.Range("H141").Value = 141091.72
$ws.Range("J141").Value = 149998.67
$ws.Range("L141").Value = 149998.67
$ws.Range("N141").Value = -160358.67
